# Auto-generated edit script: update cryptos price/volume data
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.493.02'
$ws.Range("E2").Value = '  +0.39%  '

# Row 3
$ws.Range("D3").Value = '1.852.47'
$ws.Range("E3").Value = '  -0.43%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
Set-TextValue $ws.Range("D5") '233.34'
$ws.Range("E5").Value = '  -0.65%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.001'
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4692'
$ws.Range("E7").Value = '  -0.83%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.2739'
$ws.Range("E8").Value = '  -0.37%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.06332'
$ws.Range("E9").Value = '  -1.70%  '

# Row 10
Set-TextValue $ws.Range("D10") '17.36'
$ws.Range("E10").Value = '  +5.97%  '

# Row 11
$ws.Range("D11").Value = '1.845.66'
$ws.Range("E11").Value = '  -0.93%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.07434'
$ws.Range("E12").Value = '  -0.02%  '

# Row 13
Set-TextValue $ws.Range("D13") '5.085'
$ws.Range("E13").Value = '  +1.48%  '

# Row 14
$ws.Range("E14").Value = '  -1.39%  '

# Row 15
$ws.Range("E15").Value = '  -1.83%  '

# Row 16
$ws.Range("D16").Value = '30.484.05'
$ws.Range("E16").Value = '  +0.50%  '

# Row 17
Set-TextValue $ws.Range("D17") '242.70'
$ws.Range("E17").Value = '  +4.03%  '

# Row 18
$ws.Range("E18").Value = '  +0.07%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.68'
$ws.Range("E19").Value = '  -0.76%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.000007321'
$ws.Range("E20").Value = '  -1.16%  '

# Row 21
Set-TextValue $ws.Range("D21") '1.002'
$ws.Range("E21").Value = '  +0.25%  '

# Row 22
Set-TextValue $ws.Range("D22") '4.951'
$ws.Range("E22").Value = '  -1.39%  '

# Row 23
Set-TextValue $ws.Range("D23") '5.992'
$ws.Range("E23").Value = '  -0.48%  '

# Row 24
Set-TextValue $ws.Range("D24") '9.278'
$ws.Range("E24").Value = '  -0.21%  '

# Row 25
Set-TextValue $ws.Range("D25") '163.53'
$ws.Range("E25").Value = '  -1.53%  '

# Row 26
Set-TextValue $ws.Range("D26") '18.02'
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
Set-TextValue $ws.Range("D27") '1.881'
$ws.Range("E27").Value = '  -0.80%  '

# Row 28
Set-TextValue $ws.Range("D28") '0.1009'
$ws.Range("E28").Value = '  -3.44%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.376'
$ws.Range("E29").Value = '  -0.71%  '

# Row 30
Set-TextValue $ws.Range("D30") '4.032'
$ws.Range("E30").Value = '  -2.79%  '

# Row 31
Set-TextValue $ws.Range("D31") '3.841'
$ws.Range("E31").Value = '  -2.52%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.04895'
$ws.Range("E32").Value = '  -0.51%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.139'
$ws.Range("E33").Value = '  -1.28%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.7057'
$ws.Range("E34").Value = '  -3.26%  '

# Row 35
Set-TextValue $ws.Range("D35") '2.709'
$ws.Range("E35").Value = '  +0.45%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.01902'
$ws.Range("E36").Value = '  -0.40%  '

# Row 37
Set-TextValue $ws.Range("D37") '2.683'
$ws.Range("E37").Value = '  +1.19%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.8739'
$ws.Range("E38").Value = '  -3.83%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.972'
$ws.Range("E39").Value = '  -0.50%  '

# Row 40
Set-TextValue $ws.Range("D40") '105.02'

# Row 41
$ws.Range("E41").Value = '  +0.26%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D42") '0.4064'
$ws.Range("E42").Value = '  -1.53%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D43") '5.505'
$ws.Range("E43").Value = '  -1.14%  '

# Row 44
Set-TextValue $ws.Range("D44") '7.229'
$ws.Range("E44").Value = '  +0.89%  '

# Row 45
Set-TextValue $ws.Range("D45") '62.82'
$ws.Range("E45").Value = '  +2.47%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.1201'
$ws.Range("E46").Value = '  -0.77%  '

# Row 47
Set-TextValue $ws.Range("D47") '8.599'
$ws.Range("E47").Value = '  -0.96%  '

# Row 48
Set-TextValue $ws.Range("D48") '33.34'
$ws.Range("E48").Value = '  -0.24%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.05535'
$ws.Range("E49").Value = '  -1.00%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.357'
$ws.Range("E50").Value = '  -3.84%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.3671'
$ws.Range("E51").Value = '  -1.24%  '
